# Updated symbol list on Fri Dec 16 04:58:32 UTC 2022 with GitHub Actions
#
# Refreshed "Price" (column D) quotes for several coins, and re-sorted three
# rows (41-43) so KickToken / BKEXToken / CEJI rotate position with refreshed
# prices, links and rank labels.
#
# Price cells are text (e.g. "0.1120" must stay "0.1120", not become 0.112),
# so numeric-looking values are written with a leading apostrophe - exactly
# like typing '0.1120 into a cell in Excel - to force text storage instead
# of the value being auto-coerced to a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'260.24"
$ws.Range("D3").Value = "'22.76"
$ws.Range("D4").Value = "'6.207"
$ws.Range("D5").Value = "'0.06136"
$ws.Range("D6").Value = "'6.739"
$ws.Range("D8").Value = "'1.368"
$ws.Range("D9").Value = "'0.7988"
$ws.Range("D10").Value = "'0.1578"
$ws.Range("D11").Value = "'0.08110"
$ws.Range("D12").Value = "'0.03482"
$ws.Range("D13").Value = "'0.03097"
$ws.Range("D14").Value = "'0.09331"
$ws.Range("D15").Value = "'3.865"
$ws.Range("D16").Value = "'0.001688"
$ws.Range("D17").Value = "'0.04789"
$ws.Range("D18").Value = "'0.0006152"
$ws.Range("D19").Value = "'0.006189"
$ws.Range("D20").Value = "'0.001092"
$ws.Range("D21").Value = "'0.004064"
$ws.Range("D23").Value = "'3.690"
$ws.Range("D24").Value = "'2.215"
$ws.Range("D40").Value = "'0.04615"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.007122"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1117"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.003601"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("D46").Value = "'0.00005940"
$ws.Range("D48").Value = "'0.7002"
$ws.Range("D49").Value = "'0.08974"
$ws.Range("D50").Value = "'0.00002101"
